$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3955935532374564
$ws.Range("D2").Value = 0.4321333824756292
$ws.Range("G2").Value = 0.4760219657335256
$ws.Range("H2").Value = 0.998

$ws.Range("B3").Value = 0.07514644587374561
$ws.Range("D3").Value = 0.2119198634755611
$ws.Range("G3").Value = 0.4760219657335256
$ws.Range("H3").Value = 0.998

$ws.Range("B4").Value = 0.04215534119371416
$ws.Range("D4").Value = 0.136128825357167
$ws.Range("G4").Value = 0.4760219657335256
$ws.Range("H4").Value = 0.998

$ws.Range("B5").Value = 0.07796894984218643
$ws.Range("D5").Value = 0.1911874935925044
$ws.Range("G5").Value = 0.4760219657335256
$ws.Range("H5").Value = 0.998
